# Replace the Slovenian Orion observation campaign sentence throughout
# the document (it appears identically in 4 separate paragraphs).

$d = $word.ActiveDocument

$dash = [char]0x2013

$old = "m" + [char]0x00F4 + [char]0x017E + "ete pozorova" + [char]0x0165 + " s" + [char]0x00FA + "hvezdie ozvezdje Orion 2022: 16." + $dash + "25. januar, 14." + $dash + "23. februar, 14." + $dash + "24. marec"

$new = "2022: Datumi kampanje za opazovanje ozvezdje Orion: 16." + $dash + "25. januar, 14." + $dash + "23. februar, 14." + $dash + "24. marec"

$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

while ($range.Find.Found) {
    $range.Collapse(0)
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
